$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new column order/labels -> Category, Amount, Date
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Amount"
$ws.Range("C1").Value = "Date"

# Row 2: expense entry (Food, 3000, date kept as text like the source data)
$ws.Range("A2").Value = "Food"
$ws.Range("B2").Value = 3000
$ws.Range("C2").Value = "'2025-07-21"

# Row 3: income entry (Food, 3000, date kept as text like the source data)
$ws.Range("A3").Value = "Food"
$ws.Range("B3").Value = 3000
$ws.Range("C3").Value = "'2025-07-21"
